{"js": "// Update the referenced build/changeset number from C35784 to C35789\n// (the digits run \"35784\" -> \"35789\"), keeping the existing bold\n// formatting of the run intact.\nconst body = context.document.body;\nconst buildNumMatches = body.search(\"35784\", { matchCase: true });\nbuildNumMatches.load(\"text\");\nawait context.sync();\n\nif (buildNumMatches.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for '35784', found \" + buildNumMatches.items.length\n  );\n}\nbuildNumMatches.items[0].insertText(\"35789\", \"Replace\");\nawait context.sync();\n\n// The footer of the third section (the body/content section, whose page\n// numbering restarts at 1) has a stale cached PAGE field result of \"6\".\n// Update the cached text to \"1\" to match the new pagination.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst lastSection = sections.items[sections.items.length - 1];\nconst footer = lastSection.getFooter(\"Primary\");\nconst pageNumMatches = footer.search(\"6\", { matchCase: true });\npageNumMatches.load(\"text\");\nawait context.sync();\n\nif (pageNumMatches.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for '6' in the footer, found \" + pageNumMatches.items.length\n  );\n}\npageNumMatches.items[0].insertText(\"1\", \"Replace\");\nawait context.sync();\n", "ps1": "# Update the referenced build/changeset number from C35784 to C35789\n# (the digits run \"35784\" -> \"35789\"), keeping the existing bold\n# formatting of the run intact.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"35784\"\n$find.Replacement.Text = \"35789\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, $null, $false, $null, 2) | Out-Null\n\n# The footer of the third/last section (the body/content section, whose\n# page numbering restarts at 1) has a stale cached PAGE field result of\n# \"6\". Update the cached text to \"1\" to match the new pagination.\n$lastSection = $d.Sections.Item($d.Sections.Count)\n$footer = $lastSection.Footers.Item(1)\n\n$ffind = $footer.Range.Find\n$ffind.ClearFormatting()\n$ffind.Replacement.ClearFormatting()\n$ffind.Text = \"6\"\n$ffind.Replacement.Text = \"1\"\n$ffind.Execute($null, $true, $false, $false, $false, $false, $true, $null, $false, $null, 2) | Out-Null\n"}
